# Generate Report for Handback
# Updates the handoff/handback timestamp strings recorded in the report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 9434bfbc... row (G2)
$wsOverview.Range("G2").Value = "2016-08-26 09:10:48"

# zh-cn sheet: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
# for the 9434bfbc... row
$wsZhCn.Range("H2").Value = "2016-08-26 09:10:43"
$wsZhCn.Range("K2").Value = "2016-08-26 09:11:17"

# de-de sheet: Correspond Handoff Datetime (H2) / Correspond Handback DateTime (K2)
# for the 9434bfbc... row
$wsDeDe.Range("H2").Value = "2016-08-26 09:10:48"
$wsDeDe.Range("K2").Value = "2016-08-26 09:11:24"
